# Add a new restaurant entry ("Neuse River Brewing" - burgers) to the
# bottom of the data table on Sheet1, and move the active selection to
# reflect where the author's cursor ended up (C11) after scrolling down
# while reviewing/adding the new entries (topLeftCell -> A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended right after the existing last row (27).
$ws.Range("A28").Value = "Neuse River Brewing"
$ws.Range("B28").Value = "Raleigh"
$ws.Range("C28").Value = "Neuse River Burger"
$ws.Range("D28").Value = "Brassiere/Burgers"
$ws.Range("E28").Value = 35.804566959442603
$ws.Range("F28").Value = -78.632520307935593

# Leave the selection where the author's cursor landed afterwards.
$ws.Range("C11").Select()
